$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: insert a new "2022-Q1" sheet right after "2021-Q4" (and
# therefore right before "总计"). Cloning "2021-Q4" gives us the same
# sheet formatting (borders/fonts/alignment) the other quarter sheets
# use, so we only have to overwrite the data afterwards.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newQ1 = $wb.Worksheets.Item(3)
$newQ1.Name = "2022-Q1"

# Columns B-G hold text values in this workbook (fund code, fund name,
# fund size, position, position ratio, held market value) even though
# several of them look numeric ("12.96", "0.8657", ...). Force a text
# number-format before writing so those values are NOT auto-coerced
# into numbers.
$dataRange = $newQ1.Range("B2:G6")
$dataRange.NumberFormat = "@"

$newQ1.Range("A2").Value = 0
$newQ1.Range("B2").Value = "161611"
$newQ1.Range("C2").Value = "融通内需驱动混合"
$newQ1.Range("D2").Value = "12.96"
$newQ1.Range("E2").Value = "65.58"
$newQ1.Range("F2").Value = "6.68"
$newQ1.Range("G2").Value = "0.8657"
$newQ1.Range("H2").Value = 1

$newQ1.Range("A3").Value = 1
$newQ1.Range("B3").Value = "002252"
$newQ1.Range("C3").Value = "融通成长30灵活配置混合"
$newQ1.Range("D3").Value = "1.70"
$newQ1.Range("E3").Value = "78.26"
$newQ1.Range("F3").Value = "9.54"
$newQ1.Range("G3").Value = "0.1622"
$newQ1.Range("H3").Value = 1

$newQ1.Range("A4").Value = 2
$newQ1.Range("B4").Value = "006165"
$newQ1.Range("C4").Value = "建信中证1000指数增强A"
$newQ1.Range("D4").Value = "2.75"
$newQ1.Range("E4").Value = "93.00"
$newQ1.Range("F4").Value = "1.44"
$newQ1.Range("G4").Value = "0.0396"
$newQ1.Range("H4").Value = 3

$newQ1.Range("A5").Value = 3
$newQ1.Range("B5").Value = "006166"
$newQ1.Range("C5").Value = "建信中证1000指数增强C"
$newQ1.Range("D5").Value = "0.65"
$newQ1.Range("E5").Value = "93.00"
$newQ1.Range("F5").Value = "1.44"
$newQ1.Range("G5").Value = "0.0094"
$newQ1.Range("H5").Value = 3

$newQ1.Range("A6").Value = 4
$newQ1.Range("B6").Value = "013442"
$newQ1.Range("C6").Value = "建信中证1000指数增强E"
$newQ1.Range("D6").Value = "0.02"
$newQ1.Range("E6").Value = "93.00"
$newQ1.Range("F6").Value = "1.44"
$newQ1.Range("G6").Value = "0.0003"
$newQ1.Range("H6").Value = 3

# Values are committed now; drop back to the default "Normal" style so
# no stray number-format lingers on the cells (the template columns
# carry no explicit style either).
$dataRange.Style = "Normal"

# ------------------------------------------------------------------
# Step 2: the "总计" (totals) sheet gets a new top data row for
# 2022-Q1, pushing the existing 2021-Q4 / 2021-Q3 rows down by one.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Snapshot the existing two data rows (read via Value2 - .Value's
# getter returns a COM property descriptor rather than the scalar in
# this host, Value2 is the reliable read path).
$oldB2 = $totals.Range("B2").Value2
$oldC2 = $totals.Range("C2").Value2
$oldD2 = $totals.Range("D2").Value2
$oldB3 = $totals.Range("B3").Value2
$oldC3 = $totals.Range("C3").Value2
$oldD3 = $totals.Range("D3").Value2

# Row 4 is brand new - clone row 2's formatting (bold index style "2")
# into it before writing values so it matches the other index cells.
$totals.Range("A2").Copy()
$totals.Range("A4").PasteSpecial(-4122)

# Old row 3 (2021-Q3) -> row 4
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = $oldB3
$totals.Range("C4").Value = $oldC3
$totals.Range("D4").Value = $oldD3

# Old row 2 (2021-Q4) -> row 3
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = $oldB2
$totals.Range("C3").Value = $oldC2
$totals.Range("D3").Value = $oldD2

# New row 2: 2022-Q1
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 1.08
